$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header row (A1, C1, D1, E1, J1, K1, L1) into row 12 (A12:G12)
$ws.Range("A12").Value = $ws.Range("A1").Value2
$ws.Range("B12").Value = $ws.Range("C1").Value2
$ws.Range("C12").Value = $ws.Range("D1").Value2
$ws.Range("D12").Value = $ws.Range("E1").Value2
$ws.Range("E12").Value = $ws.Range("J1").Value2
$ws.Range("F12").Value = $ws.Range("K1").Value2
$ws.Range("G12").Value = $ws.Range("L1").Value2

# Copy data rows 2-8 (A, C, D, E, J, K, L) as values into rows 13-19 (A-G)
for ($i = 0; $i -lt 7; $i++) {
    $srcRow = 2 + $i
    $dstRow = 13 + $i

    $ws.Range("A$dstRow").Value = $ws.Range("A$srcRow").Value2
    $ws.Range("B$dstRow").Value = $ws.Range("C$srcRow").Value2
    $ws.Range("C$dstRow").Value = $ws.Range("D$srcRow").Value2
    $ws.Range("D$dstRow").Value = $ws.Range("E$srcRow").Value2
    $ws.Range("E$dstRow").Value = $ws.Range("J$srcRow").Value2
    $ws.Range("F$dstRow").Value = $ws.Range("K$srcRow").Value2
    $ws.Range("G$dstRow").Value = $ws.Range("L$srcRow").Value2
}

# Update the view: select A12:G19 (A12 becomes the active cell) and scroll
# the window so row 6 is at the top (matches the saved sheetView's
# topLeftCell="A6" / selection activeCell="A12" sqref="A12:G19").
$ws.Range("A12:G19").Select()
try {
    $excel.ActiveWindow.ScrollRow = 6
} catch {
    # Older/limited hosts may not support programmatic scroll position;
    # the selection above is the functionally important part.
}
